$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shifted rows 393-484 (date + price columns) ---
$ws.Range("D393").Value = 44543
$ws.Range("J393").Value = 2500
$ws.Range("K393").Value = 600
$ws.Range("M393").Value = 650
$ws.Range("P393").Value = 650
$ws.Range("D394").Value = 44543
$ws.Range("J394").Value = 1460
$ws.Range("D395").Value = 44321
$ws.Range("J395").Value = 3400
$ws.Range("K395").Value = 650
$ws.Range("M395").Value = 675
$ws.Range("P395").Value = 675
$ws.Range("D396").Value = 44321
$ws.Range("J396").Value = 1800
$ws.Range("D397").Value = 44385
$ws.Range("J397").Value = 2200
$ws.Range("K397").Value = 600
$ws.Range("L397").Value = 700
$ws.Range("M397").Value = 650
$ws.Range("P397").Value = 650
$ws.Range("D398").Value = 44385
$ws.Range("J398").Value = 1300
$ws.Range("K398").Value = 500
$ws.Range("L398").Value = 550
$ws.Range("M398").Value = 525
$ws.Range("P398").Value = 525
$ws.Range("D399").Value = 44278
$ws.Range("J399").Value = 2600
$ws.Range("K399").Value = 700
$ws.Range("L399").Value = 800
$ws.Range("M399").Value = 750
$ws.Range("P399").Value = 750
$ws.Range("D400").Value = 44278
$ws.Range("J400").Value = 1400
$ws.Range("K400").Value = 600
$ws.Range("L400").Value = 650
$ws.Range("M400").Value = 625
$ws.Range("P400").Value = 625
$ws.Range("D401").Value = 44308
$ws.Range("J401").Value = 2400
$ws.Range("K401").Value = 650
$ws.Range("L401").Value = 700
$ws.Range("M401").Value = 675
$ws.Range("P401").Value = 675
$ws.Range("D402").Value = 44308
$ws.Range("J402").Value = 1320
$ws.Range("K402").Value = 550
$ws.Range("L402").Value = 600
$ws.Range("M402").Value = 575
$ws.Range("P402").Value = 575
$ws.Range("D403").Value = 44281
$ws.Range("J403").Value = 3400
$ws.Range("K403").Value = 700
$ws.Range("L403").Value = 750
$ws.Range("M403").Value = 725
$ws.Range("P403").Value = 725
$ws.Range("D404").Value = 44281
$ws.Range("J404").Value = 1680
$ws.Range("K404").Value = 600
$ws.Range("L404").Value = 650
$ws.Range("M404").Value = 625
$ws.Range("P404").Value = 625
$ws.Range("D405").Value = 44187
$ws.Range("J405").Value = 2400
$ws.Range("K405").Value = 550
$ws.Range("L405").Value = 600
$ws.Range("M405").Value = 575
$ws.Range("P405").Value = 575
$ws.Range("D406").Value = 44187
$ws.Range("J406").Value = 1400
$ws.Range("K406").Value = 450
$ws.Range("L406").Value = 500
$ws.Range("M406").Value = 475
$ws.Range("P406").Value = 475
$ws.Range("D407").Value = 44474
$ws.Range("J407").Value = 2200
$ws.Range("K407").Value = 650
$ws.Range("M407").Value = 675
$ws.Range("P407").Value = 675
$ws.Range("D408").Value = 44474
$ws.Range("J408").Value = 1340
$ws.Range("K408").Value = 550
$ws.Range("L408").Value = 600
$ws.Range("M408").Value = 575
$ws.Range("P408").Value = 575
$ws.Range("D409").Value = 44446
$ws.Range("J409").Value = 2000
$ws.Range("D410").Value = 44446
$ws.Range("J410").Value = 1300
$ws.Range("D411").Value = 44350
$ws.Range("J411").Value = 2400
$ws.Range("D412").Value = 44350
$ws.Range("J412").Value = 1400
$ws.Range("D413").Value = 44529
$ws.Range("J413").Value = 2500
$ws.Range("D414").Value = 44529
$ws.Range("J414").Value = 1460
$ws.Range("D415").Value = 44405
$ws.Range("J415").Value = 3400
$ws.Range("K415").Value = 600
$ws.Range("M415").Value = 650
$ws.Range("P415").Value = 650
$ws.Range("D416").Value = 44405
$ws.Range("J416").Value = 1800
$ws.Range("K416").Value = 500
$ws.Range("L416").Value = 550
$ws.Range("M416").Value = 525
$ws.Range("P416").Value = 525
$ws.Range("D417").Value = 44413
$ws.Range("J417").Value = 2000
$ws.Range("D418").Value = 44413
$ws.Range("D419").Value = 44238
$ws.Range("J419").Value = 2400
$ws.Range("K419").Value = 650
$ws.Range("L419").Value = 700
$ws.Range("M419").Value = 675
$ws.Range("P419").Value = 675
$ws.Range("D420").Value = 44238
$ws.Range("J420").Value = 1400
$ws.Range("K420").Value = 550
$ws.Range("L420").Value = 600
$ws.Range("M420").Value = 575
$ws.Range("P420").Value = 575
$ws.Range("D421").Value = 44257
$ws.Range("J421").Value = 2600
$ws.Range("K421").Value = 850
$ws.Range("L421").Value = 900
$ws.Range("M421").Value = 875
$ws.Range("P421").Value = 875
$ws.Range("D422").Value = 44257
$ws.Range("J422").Value = 1500
$ws.Range("K422").Value = 750
$ws.Range("L422").Value = 800
$ws.Range("M422").Value = 775
$ws.Range("P422").Value = 775
$ws.Range("D423").Value = 44411
$ws.Range("J423").Value = 2000
$ws.Range("K423").Value = 650
$ws.Range("L423").Value = 700
$ws.Range("M423").Value = 675
$ws.Range("P423").Value = 675
$ws.Range("D424").Value = 44411
$ws.Range("J424").Value = 1400
$ws.Range("K424").Value = 550
$ws.Range("L424").Value = 600
$ws.Range("M424").Value = 575
$ws.Range("P424").Value = 575
$ws.Range("D425").Value = 44175
$ws.Range("J425").Value = 2200
$ws.Range("D426").Value = 44175
$ws.Range("D427").Value = 44196
$ws.Range("J427").Value = 2400
$ws.Range("K427").Value = 550
$ws.Range("L427").Value = 600
$ws.Range("M427").Value = 575
$ws.Range("P427").Value = 575
$ws.Range("D428").Value = 44196
$ws.Range("J428").Value = 1500
$ws.Range("K428").Value = 450
$ws.Range("L428").Value = 500
$ws.Range("M428").Value = 475
$ws.Range("P428").Value = 475
$ws.Range("D429").Value = 44200
$ws.Range("J429").Value = 3000
$ws.Range("K429").Value = 650
$ws.Range("M429").Value = 675
$ws.Range("P429").Value = 675
$ws.Range("D430").Value = 44200
$ws.Range("J430").Value = 1600
$ws.Range("K430").Value = 550
$ws.Range("L430").Value = 600
$ws.Range("M430").Value = 575
$ws.Range("P430").Value = 575
$ws.Range("D431").Value = 44459
$ws.Range("J431").Value = 2400
$ws.Range("K431").Value = 600
$ws.Range("L431").Value = 700
$ws.Range("M431").Value = 650
$ws.Range("P431").Value = 650
$ws.Range("D432").Value = 44459
$ws.Range("J432").Value = 1400
$ws.Range("K432").Value = 500
$ws.Range("L432").Value = 550
$ws.Range("M432").Value = 525
$ws.Range("P432").Value = 525
$ws.Range("D433").Value = 44258
$ws.Range("J433").Value = 3400
$ws.Range("K433").Value = 850
$ws.Range("L433").Value = 900
$ws.Range("M433").Value = 875
$ws.Range("P433").Value = 875
$ws.Range("D434").Value = 44258
$ws.Range("J434").Value = 2000
$ws.Range("K434").Value = 750
$ws.Range("L434").Value = 800
$ws.Range("M434").Value = 775
$ws.Range("P434").Value = 775
$ws.Range("D435").Value = 44298
$ws.Range("J435").Value = 3000
$ws.Range("K435").Value = 650
$ws.Range("M435").Value = 675
$ws.Range("P435").Value = 675
$ws.Range("D436").Value = 44298
$ws.Range("J436").Value = 1480
$ws.Range("K436").Value = 550
$ws.Range("L436").Value = 600
$ws.Range("M436").Value = 575
$ws.Range("P436").Value = 575
$ws.Range("D437").Value = 44432
$ws.Range("J437").Value = 2000
$ws.Range("K437").Value = 600
$ws.Range("M437").Value = 650
$ws.Range("P437").Value = 650
$ws.Range("D438").Value = 44432
$ws.Range("J438").Value = 1400
$ws.Range("K438").Value = 500
$ws.Range("L438").Value = 550
$ws.Range("M438").Value = 525
$ws.Range("P438").Value = 525
$ws.Range("D439").Value = 44428
$ws.Range("J439").Value = 3400
$ws.Range("D440").Value = 44428
$ws.Range("J440").Value = 1600
$ws.Range("K440").Value = 550
$ws.Range("L440").Value = 600
$ws.Range("M440").Value = 575
$ws.Range("P440").Value = 575
$ws.Range("D441").Value = 44340
$ws.Range("J441").Value = 3080
$ws.Range("D442").Value = 44340
$ws.Range("J442").Value = 1400
$ws.Range("D443").Value = 44224
$ws.Range("J443").Value = 2400
$ws.Range("D444").Value = 44224
$ws.Range("J444").Value = 1540
$ws.Range("D445").Value = 44329
$ws.Range("J445").Value = 2360
$ws.Range("K445").Value = 650
$ws.Range("M445").Value = 675
$ws.Range("P445").Value = 675
$ws.Range("D446").Value = 44329
$ws.Range("J446").Value = 1360
$ws.Range("D447").Value = 44452
$ws.Range("J447").Value = 3000
$ws.Range("K447").Value = 600
$ws.Range("L447").Value = 700
$ws.Range("M447").Value = 650
$ws.Range("P447").Value = 650
$ws.Range("D448").Value = 44452
$ws.Range("J448").Value = 1400
$ws.Range("K448").Value = 500
$ws.Range("L448").Value = 550
$ws.Range("M448").Value = 525
$ws.Range("P448").Value = 525
$ws.Range("D449").Value = 44270
$ws.Range("K449").Value = 850
$ws.Range("L449").Value = 900
$ws.Range("M449").Value = 875
$ws.Range("P449").Value = 875
$ws.Range("D450").Value = 44270
$ws.Range("J450").Value = 1360
$ws.Range("K450").Value = 750
$ws.Range("L450").Value = 800
$ws.Range("M450").Value = 775
$ws.Range("P450").Value = 775
$ws.Range("D451").Value = 44195
$ws.Range("J451").Value = 2800
$ws.Range("K451").Value = 550
$ws.Range("L451").Value = 600
$ws.Range("M451").Value = 575
$ws.Range("P451").Value = 575
$ws.Range("D452").Value = 44195
$ws.Range("J452").Value = 1600
$ws.Range("K452").Value = 450
$ws.Range("L452").Value = 500
$ws.Range("M452").Value = 475
$ws.Range("P452").Value = 475
$ws.Range("D453").Value = 44473
$ws.Range("J453").Value = 2500
$ws.Range("K453").Value = 650
$ws.Range("M453").Value = 675
$ws.Range("P453").Value = 675
$ws.Range("D454").Value = 44473
$ws.Range("J454").Value = 1440
$ws.Range("K454").Value = 550
$ws.Range("L454").Value = 600
$ws.Range("M454").Value = 575
$ws.Range("P454").Value = 575
$ws.Range("D455").Value = 44398
$ws.Range("J455").Value = 3200
$ws.Range("K455").Value = 600
$ws.Range("L455").Value = 700
$ws.Range("M455").Value = 650
$ws.Range("P455").Value = 650
$ws.Range("D456").Value = 44398
$ws.Range("J456").Value = 1720
$ws.Range("K456").Value = 500
$ws.Range("L456").Value = 550
$ws.Range("M456").Value = 525
$ws.Range("P456").Value = 525
$ws.Range("D457").Value = 44536
$ws.Range("J457").Value = 2500
$ws.Range("K457").Value = 550
$ws.Range("L457").Value = 600
$ws.Range("M457").Value = 575
$ws.Range("P457").Value = 575
$ws.Range("D458").Value = 44536
$ws.Range("J458").Value = 1500
$ws.Range("K458").Value = 450
$ws.Range("L458").Value = 500
$ws.Range("M458").Value = 475
$ws.Range("P458").Value = 475
$ws.Range("D459").Value = 44302
$ws.Range("J459").Value = 3400
$ws.Range("K459").Value = 650
$ws.Range("M459").Value = 675
$ws.Range("P459").Value = 675
$ws.Range("D460").Value = 44302
$ws.Range("J460").Value = 1700
$ws.Range("K460").Value = 550
$ws.Range("L460").Value = 600
$ws.Range("M460").Value = 575
$ws.Range("P460").Value = 575
$ws.Range("D461").Value = 44511
$ws.Range("J461").Value = 2200
$ws.Range("K461").Value = 600
$ws.Range("M461").Value = 650
$ws.Range("P461").Value = 650
$ws.Range("D462").Value = 44511
$ws.Range("J462").Value = 1300
$ws.Range("K462").Value = 500
$ws.Range("L462").Value = 550
$ws.Range("M462").Value = 525
$ws.Range("P462").Value = 525
$ws.Range("D463").Value = 44239
$ws.Range("J463").Value = 3200
$ws.Range("D464").Value = 44239
$ws.Range("K464").Value = 550
$ws.Range("L464").Value = 600
$ws.Range("M464").Value = 575
$ws.Range("P464").Value = 575
$ws.Range("D465").Value = 44344
$ws.Range("J465").Value = 3400
$ws.Range("K465").Value = 650
$ws.Range("M465").Value = 675
$ws.Range("P465").Value = 675
$ws.Range("D466").Value = 44344
$ws.Range("D467").Value = 44463
$ws.Range("J467").Value = 3460
$ws.Range("D468").Value = 44463
$ws.Range("J468").Value = 1600
$ws.Range("D469").Value = 44357
$ws.Range("J469").Value = 2400
$ws.Range("D470").Value = 44357
$ws.Range("J470").Value = 1360
$ws.Range("D471").Value = 44371
$ws.Range("J471").Value = 2200
$ws.Range("D472").Value = 44371
$ws.Range("J472").Value = 1380
$ws.Range("D473").Value = 44365
$ws.Range("J473").Value = 3400
$ws.Range("D474").Value = 44365
$ws.Range("D475").Value = 44454
$ws.Range("J475").Value = 3100
$ws.Range("K475").Value = 600
$ws.Range("L475").Value = 700
$ws.Range("M475").Value = 650
$ws.Range("P475").Value = 650
$ws.Range("D476").Value = 44454
$ws.Range("J476").Value = 1600
$ws.Range("K476").Value = 500
$ws.Range("L476").Value = 550
$ws.Range("M476").Value = 525
$ws.Range("P476").Value = 525
$ws.Range("D477").Value = 44194
$ws.Range("J477").Value = 2200
$ws.Range("K477").Value = 550
$ws.Range("L477").Value = 600
$ws.Range("M477").Value = 575
$ws.Range("P477").Value = 575
$ws.Range("D478").Value = 44194
$ws.Range("J478").Value = 1400
$ws.Range("K478").Value = 450
$ws.Range("L478").Value = 500
$ws.Range("M478").Value = 475
$ws.Range("P478").Value = 475
$ws.Range("D479").Value = 44313
$ws.Range("J479").Value = 2800
$ws.Range("K479").Value = 650
$ws.Range("L479").Value = 700
$ws.Range("M479").Value = 675
$ws.Range("P479").Value = 675
$ws.Range("D480").Value = 44313
$ws.Range("J480").Value = 1500
$ws.Range("K480").Value = 550
$ws.Range("L480").Value = 600
$ws.Range("M480").Value = 575
$ws.Range("P480").Value = 575
$ws.Range("D481").Value = 44518
$ws.Range("J481").Value = 2200
$ws.Range("K481").Value = 550
$ws.Range("L481").Value = 600
$ws.Range("M481").Value = 575
$ws.Range("P481").Value = 575
$ws.Range("D482").Value = 44518
$ws.Range("J482").Value = 1340
$ws.Range("K482").Value = 450
$ws.Range("L482").Value = 500
$ws.Range("M482").Value = 475
$ws.Range("P482").Value = 475
$ws.Range("D483").Value = 44540
$ws.Range("J483").Value = 2400
$ws.Range("K483").Value = 600
$ws.Range("L483").Value = 700
$ws.Range("M483").Value = 650
$ws.Range("P483").Value = 650
$ws.Range("D484").Value = 44540
$ws.Range("J484").Value = 1200
$ws.Range("K484").Value = 500
$ws.Range("L484").Value = 550
$ws.Range("M484").Value = 525
$ws.Range("P484").Value = 525

# --- Append new rows 485-486 ---
$ws.Range("A485").Value = 8
$ws.Range("B485").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C485").Value = 'Coquimbo'
$ws.Range("D485").Value = 44272
$ws.Range("E485").Value = 4
$ws.Range("F485").Value = 100112023
$ws.Range("G485").Value = 'Brócoli'
$ws.Range("H485").Value = 'Sin especificar'
$ws.Range("I485").Value = 'Primera'
$ws.Range("J485").Value = 3400
$ws.Range("K485").Value = 850
$ws.Range("L485").Value = 900
$ws.Range("M485").Value = 875
$ws.Range("N485").Value = '$/unidad'
$ws.Range("O485").Value = 'Provincia del Elquí'
$ws.Range("P485").Value = 875
$ws.Range("Q485").Value = 1
$ws.Range("R485").Value = 'Hortaliza'
$ws.Range("D485").NumberFormat = $ws.Range("D484").NumberFormat

$ws.Range("A486").Value = 8
$ws.Range("B486").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C486").Value = 'Coquimbo'
$ws.Range("D486").Value = 44272
$ws.Range("E486").Value = 4
$ws.Range("F486").Value = 100112023
$ws.Range("G486").Value = 'Brócoli'
$ws.Range("H486").Value = 'Sin especificar'
$ws.Range("I486").Value = 'Segunda'
$ws.Range("J486").Value = 1800
$ws.Range("K486").Value = 750
$ws.Range("L486").Value = 800
$ws.Range("M486").Value = 775
$ws.Range("N486").Value = '$/unidad'
$ws.Range("O486").Value = 'Provincia del Elquí'
$ws.Range("P486").Value = 775
$ws.Range("Q486").Value = 1
$ws.Range("R486").Value = 'Hortaliza'
$ws.Range("D486").NumberFormat = $ws.Range("D484").NumberFormat

Write-Output "done"